$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4418.6313
$ws.Range("I51").Value = 1611.5555
$ws.Range("J51").Value = 6945
$ws.Range("K51").Value = 1611.5555
$ws.Range("L51").Value = 6945
$ws.Range("M51").Value = -1127.5555
$ws.Range("N51").Value = -7913
$ws.Range("H64").Value = 3660.6924
$ws.Range("I64").Value = 3198
$ws.Range("J64").Value = 3949.875
$ws.Range("K64").Value = 3198
$ws.Range("L64").Value = 3949.875
$ws.Range("M64").Value = -2950
$ws.Range("N64").Value = -4445.875
$ws.Range("H67").Value = 3660.6924
$ws.Range("I67").Value = 3198
$ws.Range("J67").Value = 3949.875
$ws.Range("K67").Value = 3198
$ws.Range("L67").Value = 3949.875
$ws.Range("M67").Value = -2340
$ws.Range("N67").Value = -5665.875
$ws.Range("H132").Value = 2031.0857
$ws.Range("I132").Value = 1060.9122
$ws.Range("J132").Value = 6284.923
$ws.Range("K132").Value = 3182.7366
$ws.Range("L132").Value = 18854.769
$ws.Range("M132").Value = -652.7366000000002
$ws.Range("N132").Value = -23914.769
$ws.Range("H137").Value = 3034.318
$ws.Range("I137").Value = 2638.353
$ws.Range("J137").Value = 4380.6
$ws.Range("K137").Value = 7915.059
$ws.Range("L137").Value = 13141.8
$ws.Range("M137").Value = -5365.059
$ws.Range("N137").Value = -18241.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1491.9333
$ws.Range("I61").Value = 1491.3572
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 1491.3572
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -1279.3572
$ws.Range("N61").Value = -1924
$ws.Range("H63").Value = 3500
$ws.Range("I63").Value = 2566.6667
$ws.Range("J63").Value = 4200
$ws.Range("K63").Value = 2566.6667
$ws.Range("L63").Value = 4200
$ws.Range("M63").Value = -1880.6667
$ws.Range("N63").Value = -5572
$ws.Range("H66").Value = 3500
$ws.Range("I66").Value = 2566.6667
$ws.Range("J66").Value = 4200
$ws.Range("K66").Value = 12833.3335
$ws.Range("L66").Value = 21000
$ws.Range("M66").Value = -9401.333500000001
$ws.Range("N66").Value = -27864
$ws.Range("H74").Value = 48991.906
$ws.Range("I74").Value = 91837.82000000001
$ws.Range("J74").Value = 1861.4
$ws.Range("K74").Value = 91837.82000000001
$ws.Range("L74").Value = 1861.4
$ws.Range("M74").Value = -90963.82000000001
$ws.Range("N74").Value = -3609.4
$ws.Range("H77").Value = 48991.906
$ws.Range("I77").Value = 91837.82000000001
$ws.Range("J77").Value = 1861.4
$ws.Range("K77").Value = 459189.1
$ws.Range("L77").Value = 9307
$ws.Range("M77").Value = -454821.1
$ws.Range("N77").Value = -18043
$ws.Range("H132").Value = 2516.75
$ws.Range("I132").Value = 2175.3572
$ws.Range("J132").Value = 3313.3333
$ws.Range("K132").Value = 6526.071599999999
$ws.Range("L132").Value = 9939.999899999999
$ws.Range("M132").Value = -3996.071599999999
$ws.Range("N132").Value = -14999.9999
$ws.Range("H136").Value = 1491.9333
$ws.Range("I136").Value = 1491.3572
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 4474.071599999999
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -1924.071599999999
$ws.Range("N136").Value = -9600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1623.8889
$ws.Range("I107").Value = 1469.4
$ws.Range("J107").Value = 1817
$ws.Range("K107").Value = 1469.4
$ws.Range("L107").Value = 1817
$ws.Range("M107").Value = 450.5999999999999
$ws.Range("N107").Value = -5657
$ws.Range("H134").Value = 3504.9443
$ws.Range("I134").Value = 3009.8914
$ws.Range("J134").Value = 6351.5
$ws.Range("K134").Value = 9029.674199999999
$ws.Range("L134").Value = 19054.5
$ws.Range("M134").Value = -6494.674199999999
$ws.Range("N134").Value = -24124.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 44119044
$ws.Range("I31").Value = 52632372
$ws.Range("J31").Value = 33335494
$ws.Range("K31").Value = 52632372
$ws.Range("L31").Value = 33335494
$ws.Range("M31").Value = -52632077
$ws.Range("N31").Value = -33336084
$ws.Range("H34").Value = 44119044
$ws.Range("I34").Value = 52632372
$ws.Range("J34").Value = 33335494
$ws.Range("K34").Value = 52632372
$ws.Range("L34").Value = 33335494
$ws.Range("M34").Value = -52632170
$ws.Range("N34").Value = -33335898
$ws.Range("H58").Value = 1220.7941
$ws.Range("I58").Value = 1179.7894
$ws.Range("J58").Value = 1272.7333
$ws.Range("K58").Value = 1179.7894
$ws.Range("L58").Value = 1272.7333
$ws.Range("M58").Value = -976.7893999999999
$ws.Range("N58").Value = -1678.7333
$ws.Range("H105").Value = 1445
$ws.Range("I105").Value = 1300
$ws.Range("J105").Value = 1493.3334
$ws.Range("K105").Value = 1300
$ws.Range("L105").Value = 1493.3334
$ws.Range("M105").Value = 447
$ws.Range("N105").Value = -4987.3334
$ws.Range("H132").Value = 2203.125
$ws.Range("I132").Value = 1922.878
$ws.Range("J132").Value = 3844.5715
$ws.Range("K132").Value = 5768.634
$ws.Range("L132").Value = 11533.7145
$ws.Range("M132").Value = -3238.634
$ws.Range("N132").Value = -16593.7145
$ws.Range("H134").Value = 1580.3214
$ws.Range("I134").Value = 1645.2273
$ws.Range("J134").Value = 1342.3334
$ws.Range("K134").Value = 4935.6819
$ws.Range("L134").Value = 4027.0002
$ws.Range("M134").Value = -2400.6819
$ws.Range("N134").Value = -9097.0002
$ws.Range("H136").Value = 1220.7941
$ws.Range("I136").Value = 1179.7894
$ws.Range("J136").Value = 1272.7333
$ws.Range("K136").Value = 3539.3682
$ws.Range("L136").Value = 3818.199900000001
$ws.Range("M136").Value = -989.3681999999999
$ws.Range("N136").Value = -8918.1999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2581.6875
$ws.Range("J109").Value = 3321.7778
$ws.Range("L109").Value = 9965.3334
$ws.Range("N109").Value = -12045.3334
$ws.Range("H131").Value = 960.83075
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 960.83075
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2882.49225
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12962.49225

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2763.6758
$ws.Range("I132").Value = 2437.2
$ws.Range("K132").Value = 7311.599999999999
$ws.Range("M132").Value = -4781.599999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5939
$ws.Range("I132").Value = 7129.143
$ws.Range("K132").Value = 21387.429
$ws.Range("M132").Value = -18857.429
$ws.Range("H136").Value = 11496041
$ws.Range("I136").Value = 25642642
$ws.Range("J136").Value = 1928.125
$ws.Range("K136").Value = 76927926
$ws.Range("L136").Value = 5784.375
$ws.Range("M136").Value = -76925376
$ws.Range("N136").Value = -10884.375
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 20934642
$ws.Range("I122").Value = 35715244
$ws.Range("J122").Value = 241802
$ws.Range("K122").Value = 107145732
$ws.Range("L122").Value = 725406
$ws.Range("M122").Value = -107143282
$ws.Range("N122").Value = -730306
$ws.Range("H126").Value = 153847700
$ws.Range("I126").Value = 125000710
$ws.Range("J126").Value = 200002880
$ws.Range("K126").Value = 375002130
$ws.Range("L126").Value = 600008640
$ws.Range("M126").Value = -374999660
$ws.Range("N126").Value = -600013580
$ws.Range("H132").Value = 3574605.2
$ws.Range("I132").Value = 4447544.5
$ws.Range("J132").Value = 3490.4546
$ws.Range("K132").Value = 13342633.5
$ws.Range("L132").Value = 10471.3638
$ws.Range("M132").Value = -13340103.5
$ws.Range("N132").Value = -15531.3638
$ws.Range("H136").Value = 27029334
$ws.Range("I136").Value = 34485376
$ws.Range("J136").Value = 1183.5
$ws.Range("K136").Value = 103456128
$ws.Range("L136").Value = 3550.5
$ws.Range("M136").Value = -103453578
$ws.Range("N136").Value = -8650.5
